$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as a new row before the
# current row 68, pushing every subsequent row down by one (old row 147
# becomes row 148). Insert the row first so formatting/shift matches what
# Excel does for Insert > Entire Row.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new observation.
$ws.Range("A68").Value = 11
$ws.Range("B68").Value = "Vega Monumental Concepción"
$ws.Range("C68").Value = "Bíobío"
$ws.Range("D68").Value = [datetime]"2022-11-30"
$ws.Range("E68").Value = 8
$ws.Range("F68").Value = 100112021
$ws.Range("G68").Value = "Ají"
$ws.Range("H68").Value = "Americana (o)"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 40
$ws.Range("K68").Value = 13000
$ws.Range("L68").Value = 14000
$ws.Range("M68").Value = 13500
$ws.Range("N68").Value = "$/caja 15 kilos"
$ws.Range("O68").Value = "Región del Maule"
$ws.Range("P68").Value = 900
$ws.Range("Q68").Value = 15
$ws.Range("R68").Value = "Hortaliza"
